# Update "Generate Report for Handback" timestamps across the three sheets.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-10-19 11:32:06"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
#              "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-10-19 11:31:54"
$wsZhCn.Range("K2").Value = "2016-10-19 11:32:36"

# de-de sheet: "Correspond Handoff Datetime" (H2) and
#              "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-10-19 11:32:06"
$wsDeDe.Range("K2").Value = "2016-10-19 11:32:53"
